# Update the "n" (count) column C with the January data.
# Every value in C2:C54 is scaled by 1.125 (i.e. multiplied by 9 and divided by 8),
# matching the before/after values described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 54
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $oldValue = $cell.Value2
    if ($null -ne $oldValue) {
        $cell.Value2 = [Math]::Round($oldValue * 1.125, 10)
    }
}
